$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.076372476977724
$ws.Range("C2").Value = 0.339357453126155
$ws.Range("D2").Value = 0.07874679504270432
$ws.Range("E2").Value = 0.1305530787407504
$ws.Range("G2").Value = 0.5634210167354468
$ws.Range("H2").Value = 0.6922552834534343
$ws.Range("L2").Value = 0.1886098154021667
$ws.Range("M2").Value = 0.2271693040413467
$ws.Range("O2").Value = 2.486489309140069
$ws.Range("B3").Value = 0.9686840919832775
$ws.Range("C3").Value = 0.3263113533540434
$ws.Range("D3").Value = 0.07140078600620825
$ws.Range("E3").Value = 0.1321566076604418
$ws.Range("G3").Value = 0.5657797926579065
$ws.Range("H3").Value = 0.6982530816225818
$ws.Range("L3").Value = 0.1859515242259775
$ws.Range("M3").Value = 0.2105733504529965
$ws.Range("O3").Value = 2.503718966320847
$ws.Range("B4").Value = 0.9025462894373959
$ws.Range("C4").Value = 0.3182910085241986
$ws.Range("D4").Value = 0.06692476485643795
$ws.Range("E4").Value = 0.133197777312301
$ws.Range("G4").Value = 0.5677494206846774
$ws.Range("H4").Value = 0.7023446836661691
$ws.Range("L4").Value = 0.1844106572097246
$ws.Range("M4").Value = 0.2004224225427507
$ws.Range("O4").Value = 2.516247831873471
$ws.Range("B5").Value = 0.8755920007286022
$ws.Range("C5").Value = 0.3150204191106241
$ws.Range("D5").Value = 0.06510940814814603
$ws.Range("E5").Value = 0.1336363064753596
$ws.Range("G5").Value = 0.5686828375597059
$ws.Range("H5").Value = 0.7041148308225615
$ws.Range("L5").Value = 0.1838057576821512
$ws.Range("M5").Value = 0.1962958991912558
$ws.Range("O5").Value = 2.521842970791425
$ws.Range("B6").Value = 0.871116152212835
$ws.Range("C6").Value = 0.3144772135390497
$ws.Range("D6").Value = 0.06480849230312913
$ws.Range("E6").Value = 0.133709984582285
$ws.Range("G6").Value = 0.5688457201001071
$ws.Range("H6").Value = 0.7044149691899122
$ws.Range("L6").Value = 0.1837067061679605
$ws.Range("M6").Value = 0.1956113080187833
$ws.Range("O6").Value = 2.52280158295828
$ws.Range("B7").Value = 0.9021827828547657
$ws.Range("C7").Value = 0.3182469088848734
$ws.Range("D7").Value = 0.06690024727247135
$ws.Range("E7").Value = 0.1332036337796231
$ws.Range("G7").Value = 0.5677614799779533
$ws.Range("H7").Value = 0.7023681403709645
$ws.Range("L7").Value = 0.1844024060651535
$ws.Range("M7").Value = 0.2003667297423348
$ws.Range("O7").Value = 2.51632130882976
$ws.Range("B8").Value = 1.039245911088415
$ws.Range("C8").Value = 0.3348613921841661
$ws.Range("D8").Value = 0.07620673315457793
$ws.Range("E8").Value = 0.1310942365537587
$ws.Range("G8").Value = 0.5641259300418824
$ws.Range("H8").Value = 0.6942384250663736
$ws.Range("L8").Value = 0.1876743084490116
$ws.Range("M8").Value = 0.2214390524380079
$ws.Range("O8").Value = 2.492024954030057
$ws.Range("B9").Value = 1.307837155645529
$ws.Range("C9").Value = 0.3673526070653566
$ws.Range("D9").Value = 0.09473140635120103
$ws.Range("E9").Value = 0.1274061754611854
$ws.Range("G9").Value = 0.5611480127608672
$ws.Range("H9").Value = 0.681542893055493
$ws.Range("L9").Value = 0.1948135789973264
$ws.Range("M9").Value = 0.2630632904758059
$ws.Range("O9").Value = 2.459886903168552
$ws.Range("B10").Value = 1.504999965907359
$ws.Range("C10").Value = 0.391157583925434
$ws.Range("D10").Value = 0.1085121102876485
$ws.Range("E10").Value = 0.124969046707626
$ws.Range("G10").Value = 0.5615124391339634
$ws.Range("H10").Value = 0.674198116314102
$ws.Range("L10").Value = 0.200498255903824
$ws.Range("M10").Value = 0.293820367744722
$ws.Range("O10").Value = 2.445781538409278
$ws.Range("B11").Value = 1.594646867712811
$ws.Range("C11").Value = 0.4019704213071975
$ws.Range("D11").Value = 0.114819117091443
$ws.Range("E11").Value = 0.1239193084341682
$ws.Range("G11").Value = 0.5622370939629064
$ws.Range("H11").Value = 0.6712880337371558
$ws.Range("L11").Value = 0.203179555385077
$ws.Range("M11").Value = 0.3078492460428137
$ws.Range("O11").Value = 2.441440406372095
$ws.Range("B12").Value = 1.62858627240206
$ws.Range("C12").Value = 0.4060623997443713
$ws.Range("D12").Value = 0.117212922501551
$ws.Range("E12").Value = 0.12353026028568
$ws.Range("G12").Value = 0.5625922241181911
$ws.Range("H12").Value = 0.6702481107028291
$ws.Range("L12").Value = 0.2042085631712069
$ws.Range("M12").Value = 0.3131667927114492
$ws.Range("O12").Value = 2.440095869657284
$ws.Range("B13").Value = 1.621277193745243
$ws.Range("C13").Value = 0.4051812394651222
$ws.Range("D13").Value = 0.1166971299991246
$ws.Range("E13").Value = 0.1236136724141718
$ws.Range("G13").Value = 0.5625121448737076
$ws.Range("H13").Value = 0.6704693153719461
$ws.Range("L13").Value = 0.2039863412548755
$ws.Range("M13").Value = 0.3120213404126559
$ws.Range("O13").Value = 2.440372111392236
$ws.Range("B14").Value = 1.597439251509968
$ws.Range("C14").Value = 0.402307124725553
$ws.Range("D14").Value = 0.1150159470977457
$ws.Range("E14").Value = 0.1238871316052961
$ws.Range("G14").Value = 0.5622646909858844
$ws.Range("H14").Value = 0.6712012344186746
$ws.Range("L14").Value = 0.2032639390240263
$ws.Range("M14").Value = 0.3082866227388976
$ws.Range("O14").Value = 2.441323785820316
$ws.Range("B15").Value = 1.582836740972027
$ws.Range("C15").Value = 0.4005462984040946
$ws.Range("D15").Value = 0.1139868875436036
$ws.Range("E15").Value = 0.1240557354004583
$ws.Range("G15").Value = 0.5621236408007491
$ws.Range("H15").Value = 0.6716576403989194
$ws.Range("L15").Value = 0.2028232240407135
$ws.Range("M15").Value = 0.3059996589656748
$ws.Range("O15").Value = 2.441945725647656
$ws.Range("B16").Value = 1.499140311478129
$ws.Range("C16").Value = 0.390450589938041
$ws.Range("D16").Value = 0.1081006992132245
$ws.Range("E16").Value = 0.1250388345652507
$ws.Range("G16").Value = 0.5614763575025137
$ws.Range("H16").Value = 0.6743969778592032
$ws.Range("L16").Value = 0.200324940495932
$ws.Range("M16").Value = 0.2929042770453094
$ws.Range("O16").Value = 2.446107082883486
$ws.Range("B17").Value = 1.44778293184828
$ws.Range("C17").Value = 0.3842528563191081
$ws.Range("D17").Value = 0.1044994711800342
$ws.Range("E17").Value = 0.1256570199817362
$ws.Range("G17").Value = 0.5612226637293816
$ws.Range("H17").Value = 0.6761879368443005
$ws.Range("L17").Value = 0.1988167019262903
$ws.Range("M17").Value = 0.2848800647117145
$ws.Range("O17").Value = 2.449192164231249
$ws.Range("B18").Value = 1.418239535432974
$ws.Range("C18").Value = 0.3806865789962615
$ws.Range("D18").Value = 0.1024317257818268
$ws.Range("E18").Value = 0.1260181305564653
$ws.Range("G18").Value = 0.5611293287755501
$ws.Range("H18").Value = 0.677258621354369
$ws.Range("L18").Value = 0.1979581770023486
$ws.Range("M18").Value = 0.2802682786410742
$ws.Range("O18").Value = 2.451161920535469
$ws.Range("B19").Value = 1.408236015288082
$ws.Range("C19").Value = 0.3794788490095584
$ws.Range("D19").Value = 0.1017322385381618
$ws.Range("E19").Value = 0.1261413492910597
$ws.Range("G19").Value = 0.5611067475425102
$ws.Range("H19").Value = 0.677628102883304
$ws.Range("L19").Value = 0.1976690379693196
$ws.Range("M19").Value = 0.2787074209081339
$ws.Range("O19").Value = 2.451862361316444
$ws.Range("B20").Value = 1.453250437777399
$ws.Range("C20").Value = 0.3849127734598596
$ws.Range("D20").Value = 0.1048824569373465
$ws.Range("E20").Value = 0.1255906390888027
$ws.Range("G20").Value = 0.5612442251538425
$ws.Range("H20").Value = 0.6759930866057857
$ws.Range("L20").Value = 0.1989763280674737
$ws.Range("M20").Value = 0.2857338924526616
$ws.Range("O20").Value = 2.448843532747929
$ws.Range("B21").Value = 1.604441261582963
$ws.Range("C21").Value = 0.40315139493066
$ws.Range("D21").Value = 0.1155096025715494
$ws.Range("E21").Value = 0.1238065803407743
$ws.Range("G21").Value = 0.562335180690738
$ws.Range("H21").Value = 0.6709845671379782
$ws.Range("L21").Value = 0.2034757557800617
$ws.Range("M21").Value = 0.3093834630393673
$ws.Range("O21").Value = 2.441036124346255
$ws.Range("B22").Value = 1.703206203907257
$ws.Range("C22").Value = 0.4150560552685647
$ws.Range("D22").Value = 0.1224869912509234
$ws.Range("E22").Value = 0.1226899313029703
$ws.Range("G22").Value = 0.5635188385152219
$ws.Range("H22").Value = 0.6680729884364212
$ws.Range("L22").Value = 0.2064959783396034
$ws.Range("M22").Value = 0.3248695222680027
$ws.Range("O22").Value = 2.437678823237206
$ws.Range("B23").Value = 1.650498348372366
$ws.Range("C23").Value = 0.4087038098715823
$ws.Range("D23").Value = 0.118760105470372
$ws.Range("E23").Value = 0.1232813963063054
$ws.Range("G23").Value = 0.5628439167860648
$ws.Range("H23").Value = 0.6695938265046806
$ws.Range("L23").Value = 0.2048767622945746
$ws.Range("M23").Value = 0.3166016867385082
$ws.Range("O23").Value = 2.439310685644955
$ws.Range("B24").Value = 1.450778632281015
$ws.Range("C24").Value = 0.38461443461145
$ws.Range("D24").Value = 0.10470930082802
$ws.Range("E24").Value = 0.1256206321113447
$ws.Range("G24").Value = 0.561234313672756
$ws.Range("H24").Value = 0.6760810505701045
$ws.Range("L24").Value = 0.1989041343558569
$ws.Range("M24").Value = 0.2853478724256249
$ws.Range("O24").Value = 2.449000538258332
$ws.Range("B25").Value = 1.235202191689552
$ws.Range("C25").Value = 0.3585736931200074
$ws.Range("D25").Value = 0.08969025500984174
$ws.Range("E25").Value = 0.128355976770324
$ws.Range("G25").Value = 0.5615069539188653
$ws.Range("H25").Value = 0.6846294369101713
$ws.Range("L25").Value = 0.1928049473751443
$ws.Range("M25").Value = 0.2517713952845497
$ws.Range("O25").Value = 2.466915473054428

Write-Host "Applied 380 kV case values"
